$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update 최종점수 (K) values
$ws.Range("K2").Value = 58.4
$ws.Range("K3").Value = 56.6
$ws.Range("K4").Value = 55.4
$ws.Range("K5").Value = 55.4

# Update MACRO_SCORE (N) values
$ws.Range("N2").Value = 54.77309453746771
$ws.Range("N3").Value = 54.77309453746771
$ws.Range("N4").Value = 54.77309453746771
$ws.Range("N5").Value = 54.77309453746771
